# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" bullet list into concise,
# impact-focused accomplishment statements, shrinking the list from six
# bullets down to four.

$d = $word.ActiveDocument

# Locate the "Impact" sub-heading under "KEY ACHIEVEMENTS AND IMPACT" and
# find the six bullet paragraphs that immediately follow it.
$count = $d.Paragraphs.Count
$implIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Impact") {
        $implIndex = $i
        break
    }
}

if ($implIndex -eq -1) {
    throw "Could not locate the 'Impact' sub-heading under KEY ACHIEVEMENTS AND IMPACT"
}

$bullet1 = $d.Paragraphs.Item($implIndex + 1)
$bullet2 = $d.Paragraphs.Item($implIndex + 2)
$bullet3 = $d.Paragraphs.Item($implIndex + 3)
$bullet4 = $d.Paragraphs.Item($implIndex + 4)
$bullet5 = $d.Paragraphs.Item($implIndex + 5)
$bullet6 = $d.Paragraphs.Item($implIndex + 6)

# Replace the text of the four bullets that survive, doing this before any
# deletion so paragraph references stay valid.
$bullet1.Range.Text = "• Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard"
$bullet2.Range.Text = "• Reduced polling margins from ±4.2% to ±2.1%"
$bullet3.Range.Text = "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"
$bullet6.Range.Text = "• Reduced polling costs while increasing quality"

# Delete bullets 4 and 5 entirely ("Built redistricting platform..." and
# "Developed longitudinal data analysis methods...") by removing the
# combined range that spans both paragraphs. Re-fetch them fresh (rather
# than reuse the original references) to be safe.
$bullet4 = $d.Paragraphs.Item($implIndex + 4)
$bullet5 = $d.Paragraphs.Item($implIndex + 5)
$deleteRange = $d.Range($bullet4.Range.Start, $bullet5.Range.End)
$deleteRange.Delete()

Write-Output "Key Achievements section updated"
